# B6-PowerPoint.pptx edit — Thu, Aug 06, 2020 12:05:29 AM
#
# 1) Three tables (slides 14, 15, 16) switch their table style from the
#    deck's single custom style ("Table_0",
#    {BEAE1C6D-6CE9-4744-860C-CFD7D6BA5E34}) to the built-in
#    "No Style, No Grid" table style ({3D504A21-3887-4A16-B0CE-F6B347593171}).
# 2) The deck's theme ("Integral" / "Red Violet" palette) is swapped for the
#    plain default "Office" color palette (the palette that used to live
#    only on the otherwise-unused notes-master theme part).

$p = $ppt.ActivePresentation

# --- 1) Re-style the three tables -----------------------------------------
$noStyleNoGridId = "{3D504A21-3887-4A16-B0CE-F6B347593171}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($noStyleNoGridId)
        }
    }
}

# --- 2) Swap the presentation's theme colors for the plain Office palette --
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

# Order matches the standard DrawingML clrScheme slot order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeRgb = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeRgb[$i - 1]
}
